$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Factory")

$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 0
